# Generate Report for Handoff
# Updates the localization status report: marks the tracked item as
# "Ready for handoff" (was "In Translation") and refreshes the related
# handoff timestamps on all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
# zh-cn / de-de status columns + "Latest HO Xliff Generate Date"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-16 08:39:17"

# --- zh-cn sheet ---
# Status + "Latest Handoff Datetime"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-16 08:39:12"

# --- de-de sheet ---
# Status + "Latest Handoff Datetime"
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-16 08:39:17"

# The longer "Ready for handoff" text no longer fits the previous column
# width, so the Status columns grow to accommodate it (as Excel's
# autofit would do after the text change).
$wsOverview.Range("E:F").ColumnWidth = 16.3333333333333
$wsZhCn.Range("C:C").ColumnWidth = 16.3333333333333
$wsDeDe.Range("C:C").ColumnWidth = 16.3333333333333
